$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix existing material mismatches (E4 and E8 swap)
$ws.Range("E4").Value = "Metal"
$ws.Range("E8").Value = "Wood"

# Apply the centered style (same as used elsewhere, s="1") to rows 35-39
$ws.Range("A2").Copy()
$ws.Range("A35:E39").PasteSpecial(-4122)

# Add new rows 40-43
$ws.Range("A40").Value = 43
$ws.Range("B40").Value = 0.5
$ws.Range("C40").Value = 2
$ws.Range("D40").Value = 5.222
$ws.Range("E40").Value = "Plastic"

$ws.Range("A41").Value = "423FF"
$ws.Range("B41").Value = 3.2
$ws.Range("C41").Value = 33
$ws.Range("D41").Value = 3
$ws.Range("E41").Value = "Plastic"

$ws.Range("A42").Value = "187A"
$ws.Range("B42").Value = 15.8
$ws.Range("C42").Value = 10
$ws.Range("D42").Value = 30
$ws.Range("E42").Value = "Wood"

$ws.Range("A43").Value = "52524A"
$ws.Range("B43").Value = 250
$ws.Range("C43").Value = 250
$ws.Range("D43").Value = 250
$ws.Range("E43").Value = "Plastic"

# Apply style to new rows too
$ws.Range("A2").Copy()
$ws.Range("A40:E43").PasteSpecial(-4122)

# Column widths
$ws.Columns.Item(1).ColumnWidth = 12.3984375
$ws.Columns.Item(5).ColumnWidth = 13

# Selection
$ws.Range("G41").Select() | Out-Null
